$wb = $excel.ActiveWorkbook

# Duplicate the "Croatia" sheet to create the new "Greece" sheet, placing it
# after the last existing sheet (Croatia is currently last).
$croatia = $wb.Worksheets.Item("Croatia")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Copy($null, $lastSheet)

$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Fill in the Greece-specific data (ticket reference + market name). The
# ticket reference is entered first so it lands before the market name in
# the shared-strings table, matching the order in which the two new
# strings were appended upstream.
$greece.Range("B4").Value = "NGC-4119/T3202"
$greece.Range("B2").Value = "Greece Market"

# Restore Croatia's selection/view state.
$croatia.Activate()
$croatia.Cells.Select()

# Make Greece the active sheet/tab, with its own selection.
$greece.Activate()
$greece.Range("E16").Select()
